$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (SKU moves from B to C, etc.)
$ws.Columns("B").EntireColumn.Insert()

# New header cell for the inserted column. Insert() already carries the
# surrounding "Heading 1" formatting into the new column, so just set the
# text.
$ws.Range("B1").Value = "Category"

# Restore the current selection seen in the saved file.
$ws.Range("B7").Select()
